# Applies the diff: a new price-report data row is inserted as row 13,
# pushing the existing rows 13:99 down to 14:100 (dimension grows from
# A1:T99 to A1:T100). The new row carries the same Mercado/Producto
# context as its neighbours plus its own price-observation data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 13:99 down to 14:100, inheriting formatting (incl. the
# date-formatted D column) from the row being pushed down.
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new observation.
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Vega Monumental Concepción"
$ws.Range("C13").Value = "Bíobío"
$ws.Range("D13").Value = 44575
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100109
$ws.Range("H13").Value = "Uva"
$ws.Range("I13").Value = 100109001
$ws.Range("J13").Value = "Uva"
$ws.Range("K13").Value = "Superior Seedless"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 10000
$ws.Range("O13").Value = 11000
$ws.Range("P13").Value = 10500
$ws.Range("Q13").Value = "$/caja 10 kilos"
$ws.Range("R13").Value = "Provincia de Limarí"
$ws.Range("S13").Value = 1050
$ws.Range("T13").Value = 10
